$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 11.146846
$ws.Range("N2").Value = 33.440538
$ws.Range("O2").Value = 0.2594806085672136
$ws.Range("P2").Value = 0.2594806085672136
$ws.Range("Q2").Value = 25.7291499372
$ws.Range("R2").Value = 231.5623494348
$ws.Range("S2").Value = 0.1722507171045196
$ws.Range("T2").Value = 0.1722507171045196

$ws.Range("O3").Value = 0.6444737471070977
$ws.Range("P3").Value = 0.6444737471070977
$ws.Range("S3").Value = 0.4278202741515426
$ws.Range("T3").Value = 0.4278202741515426

$ws.Range("O4").Value = 0.09604564432568881
$ws.Range("P4").Value = 0.09604564432568881
$ws.Range("S4").Value = 0.06375787077584312
$ws.Range("T4").Value = 0.06375787077584312

$ws.Range("M5").Value = 11.146846
$ws.Range("N5").Value = 33.440538
$ws.Range("O5").Value = 0.2594806085672136
$ws.Range("P5").Value = 0.2594806085672136
$ws.Range("Q5").Value = 13.029559436246
$ws.Range("R5").Value = 117.266034926214
$ws.Range("S5").Value = 0.08722989146269393
$ws.Range("T5").Value = 0.08722989146269391

$ws.Range("O6").Value = 0.6444737471070977
$ws.Range("P6").Value = 0.6444737471070977
$ws.Range("S6").Value = 0.2166534729555551
$ws.Range("T6").Value = 0.2166534729555551

$ws.Range("O7").Value = 0.09604564432568881
$ws.Range("P7").Value = 0.09604564432568881
$ws.Range("S7").Value = 0.03228777354984569
$ws.Range("T7").Value = 0.03228777354984568
